$wb = $excel.ActiveWorkbook

# The "Another" sheet holds headered sample data in A1:D4 (row 1 = headers,
# rows 2-4 = data). Add a new sheet "NoHeaders" directly after it, holding
# the same data but without the header row (so it ends up in A1:D3).
$src = $wb.Worksheets.Item("Another")
$new = $wb.Worksheets.Add($null, $src)
$new.Name = "NoHeaders"

for ($r = 0; $r -lt 3; $r++) {
    for ($c = 1; $c -le 4; $c++) {
        $srcCell = $src.Cells.Item($r + 2, $c)
        $dstCell = $new.Cells.Item($r + 1, $c)
        $dstCell.Value = $srcCell.Value2
    }
}

# Column D held dates formatted on "Another" - carry the same number format
# over to the copied column D on the new sheet.
$new.Range("D1:D3").NumberFormat = $src.Range("D2:D4").NumberFormat

# Update the selection on "Another" to cover the data rows too.
$src.Activate() | Out-Null
$src.Range("A2:D4").Select() | Out-Null

# Select the full data range on the new sheet, and leave it as the active tab.
$new.Activate() | Out-Null
$new.Range("A1:D3").Select() | Out-Null
